# Applies numeric corrections to the Leve profit-tracking sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Generated from the authoritative cell-level diff: updates currentAveragePrice* / LevePrice* / LeveProfit*
# columns (H-N) for the affected Leve rows on each sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4238
$ws.Range("J17").Value = 3696.8
$ws.Range("L17").Value = 11090.4
$ws.Range("N17").Value = -11426.4
$ws.Range("H98").Value = 2877.1
$ws.Range("I98").Value = 2483.125
$ws.Range("J98").Value = 4453
$ws.Range("K98").Value = 2483.125
$ws.Range("L98").Value = 4453
$ws.Range("M98").Value = -985.125
$ws.Range("N98").Value = -7449
$ws.Range("H112").Value = 1653.9385
$ws.Range("J112").Value = 1653.9385
$ws.Range("L112").Value = 4961.8155
$ws.Range("N112").Value = -7177.8155
$ws.Range("H122").Value = 2877.1
$ws.Range("I122").Value = 2483.125
$ws.Range("J122").Value = 4453
$ws.Range("K122").Value = 7449.375
$ws.Range("L122").Value = 13359
$ws.Range("M122").Value = -4999.375
$ws.Range("N122").Value = -18259
$ws.Range("H132").Value = 987.6667
$ws.Range("I132").Value = 899.65
$ws.Range("K132").Value = 2698.95
$ws.Range("M132").Value = -168.9499999999998
$ws.Range("H137").Value = 26631.795
$ws.Range("I137").Value = 733.2273
$ws.Range("K137").Value = 2199.6819
$ws.Range("M137").Value = 350.3181
$ws.Range("H139").Value = 74000
$ws.Range("J139").Value = 74000
$ws.Range("L139").Value = 74000
$ws.Range("N139").Value = -84280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4140.94
$ws.Range("I32").Value = 3271.2307
$ws.Range("J32").Value = 12934.667
$ws.Range("K32").Value = 3271.2307
$ws.Range("L32").Value = 12934.667
$ws.Range("M32").Value = -2984.2307
$ws.Range("N32").Value = -13508.667
$ws.Range("H122").Value = 2408.5293
$ws.Range("I122").Value = 1352.3572
$ws.Range("K122").Value = 4057.0716
$ws.Range("M122").Value = -1607.0716
$ws.Range("H132").Value = 1845.2245
$ws.Range("I132").Value = 1386.3846
$ws.Range("J132").Value = 2363.913
$ws.Range("K132").Value = 4159.1538
$ws.Range("L132").Value = 7091.739
$ws.Range("M132").Value = -1629.1538
$ws.Range("N132").Value = -12151.739

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 7849
$ws.Range("J80").Value = 7849
$ws.Range("L80").Value = 7849
$ws.Range("N80").Value = -9845
$ws.Range("H83").Value = 7849
$ws.Range("J83").Value = 7849
$ws.Range("L83").Value = 39245
$ws.Range("N83").Value = -49229
$ws.Range("H86").Value = 402201.4
$ws.Range("I86").Value = 2333.3333
$ws.Range("J86").Value = 1002003.5
$ws.Range("K86").Value = 2333.3333
$ws.Range("L86").Value = 1002003.5
$ws.Range("M86").Value = -1210.3333
$ws.Range("N86").Value = -1004249.5
$ws.Range("H89").Value = 402201.4
$ws.Range("I89").Value = 2333.3333
$ws.Range("J89").Value = 1002003.5
$ws.Range("K89").Value = 11666.6665
$ws.Range("L89").Value = 5010017.5
$ws.Range("M89").Value = -6050.666499999999
$ws.Range("N89").Value = -5021249.5
$ws.Range("H105").Value = 2529.28
$ws.Range("J105").Value = 3177.75
$ws.Range("L105").Value = 3177.75
$ws.Range("N105").Value = -6671.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1497.6666
$ws.Range("I132").Value = 1122.4762
$ws.Range("K132").Value = 3367.4286
$ws.Range("M132").Value = -837.4286000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 121.36842
$ws.Range("I12").Value = 70.8
$ws.Range("J12").Value = 139.42857
$ws.Range("K12").Value = 212.4
$ws.Range("L12").Value = 418.28571
$ws.Range("M12").Value = -39.39999999999998
$ws.Range("N12").Value = -764.28571
$ws.Range("H22").Value = 3625
$ws.Range("I22").Value = 3000
$ws.Range("J22").Value = 3833.3333
$ws.Range("K22").Value = 9000
$ws.Range("L22").Value = 11499.9999
$ws.Range("M22").Value = -8831
$ws.Range("N22").Value = -11837.9999
$ws.Range("H27").Value = 3625
$ws.Range("I27").Value = 3000
$ws.Range("J27").Value = 3833.3333
$ws.Range("K27").Value = 9000
$ws.Range("L27").Value = 11499.9999
$ws.Range("M27").Value = -8898
$ws.Range("N27").Value = -11703.9999
$ws.Range("H29").Value = 196
$ws.Range("I29").Value = 200
$ws.Range("J29").Value = 192
$ws.Range("K29").Value = 600
$ws.Range("L29").Value = 576
$ws.Range("M29").Value = -323
$ws.Range("N29").Value = -1130
$ws.Range("H51").Value = 1624.75
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H64").Value = 2549.75
$ws.Range("I64").Value = 1599.5
$ws.Range("K64").Value = 4798.5
$ws.Range("M64").Value = -4528.5
$ws.Range("H67").Value = 2549.75
$ws.Range("I67").Value = 1599.5
$ws.Range("K67").Value = 4798.5
$ws.Range("M67").Value = -3862.5
$ws.Range("H81").Value = 2702.5
$ws.Range("I81").Value = 1200
$ws.Range("J81").Value = 3003
$ws.Range("K81").Value = 3600
$ws.Range("L81").Value = 9009
$ws.Range("M81").Value = -2477
$ws.Range("N81").Value = -11255
$ws.Range("H84").Value = 2702.5
$ws.Range("I84").Value = 1200
$ws.Range("J84").Value = 3003
$ws.Range("K84").Value = 10800
$ws.Range("L84").Value = 27027
$ws.Range("M84").Value = -5184
$ws.Range("N84").Value = -38259
$ws.Range("H88").Value = 5116.5
$ws.Range("J88").Value = 5539.8
$ws.Range("L88").Value = 16619.4
$ws.Range("N88").Value = -17475.4
$ws.Range("H91").Value = 5116.5
$ws.Range("J91").Value = 5539.8
$ws.Range("L91").Value = 16619.4
$ws.Range("N91").Value = -19583.4
$ws.Range("H92").Value = 462.25
$ws.Range("I92").Value = 399.5
$ws.Range("K92").Value = 1198.5
$ws.Range("M92").Value = 49.5
$ws.Range("H107").Value = 501.6154
$ws.Range("I107").Value = 404.57144
$ws.Range("J107").Value = 614.8333
$ws.Range("K107").Value = 1213.71432
$ws.Range("L107").Value = 1844.4999
$ws.Range("M107").Value = 706.28568
$ws.Range("N107").Value = -5684.4999
$ws.Range("H122").Value = 1039.1538
$ws.Range("J122").Value = 1092.2174
$ws.Range("L122").Value = 9829.9566
$ws.Range("N122").Value = -14729.9566
$ws.Range("H128").Value = 240674.75
$ws.Range("I128").Value = 240674.75
$ws.Range("K128").Value = 722024.25
$ws.Range("M128").Value = -717044.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4514.8237
$ws.Range("I70").Value = 4429.5
$ws.Range("K70").Value = 4429.5
$ws.Range("M70").Value = -4159.5
$ws.Range("H73").Value = 4514.8237
$ws.Range("I73").Value = 4429.5
$ws.Range("K73").Value = 4429.5
$ws.Range("M73").Value = -3493.5
$ws.Range("H102").Value = 3249.5
$ws.Range("I102").Value = 5499.5
$ws.Range("K102").Value = 5499.5
$ws.Range("M102").Value = -3877.5
$ws.Range("H126").Value = 1826673.9
$ws.Range("I126").Value = 2060040.4
$ws.Range("J126").Value = 251449.75
$ws.Range("K126").Value = 6180121.199999999
$ws.Range("L126").Value = 754349.25
$ws.Range("M126").Value = -6177651.199999999
$ws.Range("N126").Value = -759289.25
$ws.Range("H132").Value = 1675445
$ws.Range("I132").Value = 3207837.2
$ws.Range("J132").Value = 3744.4546
$ws.Range("K132").Value = 9623511.600000001
$ws.Range("L132").Value = 11233.3638
$ws.Range("M132").Value = -9620981.600000001
$ws.Range("N132").Value = -16293.3638
$ws.Range("H136").Value = 7906.2173
$ws.Range("J136").Value = 7906.2173
$ws.Range("L136").Value = 23718.6519
$ws.Range("N136").Value = -28818.6519

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3017.75
$ws.Range("J7").Value = 3281.8
$ws.Range("L7").Value = 3281.8
$ws.Range("N7").Value = -3505.8
$ws.Range("H122").Value = 2886.7144
$ws.Range("I122").Value = 2868
$ws.Range("J122").Value = 2999
$ws.Range("K122").Value = 8604
$ws.Range("L122").Value = 8997
$ws.Range("M122").Value = -6154
$ws.Range("N122").Value = -13897
$ws.Range("H126").Value = 3017.75
$ws.Range("J126").Value = 3281.8
$ws.Range("L126").Value = 9845.400000000001
$ws.Range("N126").Value = -14785.4
$ws.Range("H132").Value = 4807.885
$ws.Range("I132").Value = 4385.3335
$ws.Range("J132").Value = 5031.5884
$ws.Range("K132").Value = 13156.0005
$ws.Range("L132").Value = 15094.7652
$ws.Range("M132").Value = -10626.0005
$ws.Range("N132").Value = -20154.7652

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4299.909
$ws.Range("I126").Value = 4841.654
$ws.Range("J126").Value = 2287.7144
$ws.Range("K126").Value = 14524.962
$ws.Range("L126").Value = 6863.1432
$ws.Range("M126").Value = -12054.962
$ws.Range("N126").Value = -11803.1432
$ws.Range("H132").Value = 1333.7142
$ws.Range("I132").Value = 1299.7188
$ws.Range("K132").Value = 3899.1564
$ws.Range("M132").Value = -1369.1564

